$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the data in columns B:G down by one row (row N+1 gets the old values of row N),
# for rows 2..11 (new row 2 gets freshly computed values below).
# We must read all old values first before overwriting, since writing would
# otherwise clobber values we still need to read.

$oldValues = @{}
for ($r = 2; $r -le 11; $r++) {
    $oldValues[$r] = @(
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2,
        $ws.Cells.Item($r, 7).Value2
    )
}

# Shift rows 10 -> 11, 9 -> 10, ..., 2 -> 3
for ($r = 11; $r -ge 3; $r--) {
    $src = $oldValues[$r - 1]
    $ws.Cells.Item($r, 2).Value2 = $src[0]
    $ws.Cells.Item($r, 3).Value2 = $src[1]
    $ws.Cells.Item($r, 4).Value2 = $src[2]
    $ws.Cells.Item($r, 5).Value2 = $src[3]
    $ws.Cells.Item($r, 6).Value2 = $src[4]
    $ws.Cells.Item($r, 7).Value2 = $src[5]
}

# New values for row 2
$ws.Cells.Item(2, 2).Value2 = 0.02072117565895826
$ws.Cells.Item(2, 3).Value2 = 0.5679342762134251
$ws.Cells.Item(2, 4).Value2 = 0.6583304034161481
$ws.Cells.Item(2, 5).Value2 = 0.8113756241200176
$ws.Cells.Item(2, 6).Value2 = 0.8333373229775733
$ws.Cells.Item(2, 7).Value2 = 19
